# Actualización automática hashcode vie abr 10 01:29:55 CEST 2020
# Update the "hashcode" values (column B) for the rows whose id (column A)
# matches the entries below. Each row is addressed explicitly by its row
# number to avoid any ambiguity from duplicate values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 89;  Value = "68439b0181d7876541c13259662e62d3" },
    @{ Row = 99;  Value = "0dc061740719f94d60c3f2fb1a76b472" },
    @{ Row = 110; Value = "aa9b18f3904e71ef4dadf88111858b4d" },
    @{ Row = 154; Value = "d09f757fce10d83c04b40c0872fa2f49" },
    @{ Row = 159; Value = "567cf77756c9ad1d2efe5d8d378938af" },
    @{ Row = 160; Value = "e6c7c00b8b6510a2b39353b93e9900e0" },
    @{ Row = 162; Value = "f26197f222ebf1ddb0efdcaf398412ff" },
    @{ Row = 222; Value = "080d67ce5b06284f1eb7f4c2be969aae" },
    @{ Row = 229; Value = "99a87e6a63c15449ba6dc97361bdc99f" },
    @{ Row = 246; Value = "bfa945bf51564579f4f6f69446a4f8b9" },
    @{ Row = 276; Value = "f5aed2cdf9c8e1604b697fc84a146194" },
    @{ Row = 278; Value = "f628fff06e904e2e47120b72d229abd7" },
    @{ Row = 281; Value = "01fe8f04b41cfb18b35169e6f343d375" },
    @{ Row = 335; Value = "11352530e667e1d92b0f0b73ab121312" },
    @{ Row = 339; Value = "95dbda5d9a8b6ad8dfae2c4599d555fd" },
    @{ Row = 411; Value = "5016892bf179990dc395a7ae5ebe6a6f" },
    @{ Row = 420; Value = "bf3569543f5afe0bd329968445d710df" },
    @{ Row = 448; Value = "c5a9f86f2d2c71529bd5c6e6dee7b713" },
    @{ Row = 523; Value = "d955a48620bab04169b9a56444fc3ee6" },
    @{ Row = 542; Value = "8077e4472b59b3f779e1f9186c1e6d00" },
    @{ Row = 561; Value = "5860f193406589f9e74d7a8a5210a51e" },
    @{ Row = 574; Value = "dbf505ff70abf2f98d208c3f09f0e5e2" },
    @{ Row = 592; Value = "320e9ebd681ed0347b22b3f3e81e84de" },
    @{ Row = 764; Value = "b62c6f676d044fe150c53af2473ab411" },
    @{ Row = 769; Value = "7eb6480a07fb4eb0723eb9269c759c30" },
    @{ Row = 776; Value = "a674c1abc8131bd1104e7863c9f31dd5" },
    @{ Row = 794; Value = "81a0b2c99fbc2c00faea8501d1da6b2a" },
    @{ Row = 819; Value = "b6c09b428d120017c20a693b7eebd5f4" },
    @{ Row = 823; Value = "3ea8f816521a69783b352cf8f251e194" },
    @{ Row = 827; Value = "8984ed957ef45588ab2b7e250414079d" },
    @{ Row = 833; Value = "3b90ab400a44cba436858271a190263b" },
    @{ Row = 835; Value = "462b9661f05db7b33cc099b42a4fe747" },
    @{ Row = 863; Value = "1cc9163b0df7e540f7f31945be6fda33" },
    @{ Row = 877; Value = "6cb1529c3773a8074f6a97dc67c2a11e" },
    @{ Row = 913; Value = "e0b86bff7849e4f28182a83327dbfef5" },
    @{ Row = 937; Value = "15b108db17024df38c360f0ccf7ecc84" }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 2).Value = $u.Value
}
